$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 677; this shifts the existing rows
# 677-771 down to 679-773, preserving all of their values.
$ws.Rows("677:678").Insert()

# Populate the two newly inserted rows with the new weekly records.
$ws.Range("A677").Value = 8
$ws.Range("B677").Value = "Terminal La Palmera de La Serena"
$ws.Range("C677").Value = "Coquimbo"
$ws.Range("D677").Value = 44776
$ws.Range("E677").Value = 4
$ws.Range("F677").Value = 100112008
$ws.Range("G677").Value = "Coliflor"
$ws.Range("H677").Value = "Sin especificar"
$ws.Range("I677").Value = "Primera"
$ws.Range("J677").Value = 2460
$ws.Range("K677").Value = 800
$ws.Range("L677").Value = 900
$ws.Range("M677").Value = 850
$ws.Range("N677").Value = "`$/unidad"
$ws.Range("O677").Value = "Provincia del Elquí"
$ws.Range("P677").Value = 850
$ws.Range("Q677").Value = 1
$ws.Range("R677").Value = "Hortaliza"

$ws.Range("A678").Value = 8
$ws.Range("B678").Value = "Terminal La Palmera de La Serena"
$ws.Range("C678").Value = "Coquimbo"
$ws.Range("D678").Value = 44776
$ws.Range("E678").Value = 4
$ws.Range("F678").Value = 100112008
$ws.Range("G678").Value = "Coliflor"
$ws.Range("H678").Value = "Sin especificar"
$ws.Range("I678").Value = "Segunda"
$ws.Range("J678").Value = 1540
$ws.Range("K678").Value = 700
$ws.Range("L678").Value = 750
$ws.Range("M678").Value = 725
$ws.Range("N678").Value = "`$/unidad"
$ws.Range("O678").Value = "Provincia del Elquí"
$ws.Range("P678").Value = 725
$ws.Range("Q678").Value = 1
$ws.Range("R678").Value = "Hortaliza"
